$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily to prevent Excel from
# auto-converting values like "1.009" / "0.7190" into numbers (which
# would strip significant trailing zeros / change the stored type).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.101.55"
$ws.Range("E2").Value = "  -1.65%  "

# Row 3
$ws.Range("D3").Value = "1.798.65"
$ws.Range("E3").Value = "  -2.26%  "

# Row 4
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.55%  "

# Row 5
$ws.Range("D5").Value = "1.009"
$ws.Range("E5").Value = "  +0.57%  "

# Row 6
$ws.Range("D6").Value = "308.06"
$ws.Range("E6").Value = "  -1.78%  "

# Row 7
$ws.Range("D7").Value = "0.4163"
$ws.Range("E7").Value = "  -1.99%  "

# Row 8
$ws.Range("D8").Value = "0.3551"
$ws.Range("E8").Value = "  -3.13%  "

# Row 9
$ws.Range("D9").Value = "0.07024"
$ws.Range("E9").Value = "  -2.96%  "

# Row 10
$ws.Range("D10").Value = "0.8413"
$ws.Range("E10").Value = "  -3.26%  "

# Row 11
$ws.Range("D11").Value = "1.930.64"
$ws.Range("E11").Value = "  +2.87%  "

# Row 12
$ws.Range("D12").Value = "20.15"
$ws.Range("E12").Value = "  -3.09%  "

# Row 13
$ws.Range("D13").Value = "5.257"
$ws.Range("E13").Value = "  -2.60%  "

# Row 14
$ws.Range("D14").Value = "6.322"
$ws.Range("E14").Value = "  -3.07%  "

# Row 15
$ws.Range("D15").Value = "0.06823"
$ws.Range("E15").Value = "  -1.58%  "

# Row 16
$ws.Range("D16").Value = "1.010"
$ws.Range("E16").Value = "  +0.56%  "

# Row 17
$ws.Range("D17").Value = "79.79"
$ws.Range("E17").Value = "  -0.42%  "

# Row 18
$ws.Range("D18").Value = "0.000008688"
$ws.Range("E18").Value = "  -3.51%  "

# Row 19
$ws.Range("D19").Value = "1.010"
$ws.Range("E19").Value = "  +0.77%  "

# Row 20
$ws.Range("D20").Value = "15.05"
$ws.Range("E20").Value = "  -2.82%  "

# Row 21
$ws.Range("D21").Value = "27.667.63"
$ws.Range("E21").Value = "  +0.03%  "

# Row 22
$ws.Range("D22").Value = "5.044"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23
$ws.Range("D23").Value = "10.72"
$ws.Range("E23").Value = "  -1.26%  "

# Row 24
$ws.Range("D24").Value = "2.108.79"
$ws.Range("E24").Value = "  -0.47%  "

# Row 25
$ws.Range("D25").Value = "1.949"
$ws.Range("E25").Value = "  -0.26%  "

# Row 26
$ws.Range("D26").Value = "152.95"
$ws.Range("E26").Value = "  -0.70%  "

# Row 27
$ws.Range("D27").Value = "18.14"
$ws.Range("E27").Value = "  -1.25%  "

# Row 28
$ws.Range("D28").Value = "5.014"
$ws.Range("E28").Value = "  -4.68%  "

# Row 29
$ws.Range("D29").Value = "112.58"
$ws.Range("E29").Value = "  -2.37%  "

# Row 30
$ws.Range("D30").Value = "1.649"
$ws.Range("E30").Value = "  -10.48%  "

# Row 31
$ws.Range("D31").Value = "0.08854"
$ws.Range("E31").Value = "  -0.27%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.7190"
$ws.Range("E32").Value = "  -7.19%  "

# Row 33
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "2.867"
$ws.Range("E33").Value = "  -2.96%  "

# Row 34
$ws.Range("D34").Value = "4.334"
$ws.Range("E34").Value = "  -5.05%  "

# Row 35
$ws.Range("D35").Value = "1.008"
$ws.Range("E35").Value = "  +0.59%  "

# Row 36
$ws.Range("D36").Value = "1.076"
$ws.Range("E36").Value = "  -6.53%  "

# Row 37
$ws.Range("D37").Value = "1.076"
$ws.Range("E37").Value = "  -2.23%  "

# Row 38
$ws.Range("D38").Value = "0.01890"
$ws.Range("E38").Value = "  -2.93%  "

# Row 39
$ws.Range("D39").Value = "0.05087"
$ws.Range("E39").Value = "  -5.46%  "

# Row 40
$ws.Range("D40").Value = "0.4925"
$ws.Range("E40").Value = "  -3.90%  "

# Row 41
$ws.Range("D41").Value = "0.1612"
$ws.Range("E41").Value = "  -2.92%  "

# Row 42
$ws.Range("D42").Value = "2.613"
$ws.Range("E42").Value = "  -7.78%  "

# Row 43
$ws.Range("D43").Value = "6.173"
$ws.Range("E43").Value = "  -8.84%  "

# Row 44
$ws.Range("D44").Value = "8.051"
$ws.Range("E44").Value = "  -5.32%  "

# Row 45
$ws.Range("D45").Value = "1.008"
$ws.Range("E45").Value = "  +0.55%  "

# Row 46
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "104.16"
$ws.Range("E46").Value = "  -1.81%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.13"
$ws.Range("E47").Value = "  -3.33%  "

# Row 48
$ws.Range("D48").Value = "0.06309"
$ws.Range("E48").Value = "  -3.48%  "

# Row 49
$ws.Range("D49").Value = "0.4519"
$ws.Range("E49").Value = "  -4.11%  "

# Row 50
$ws.Range("D50").Value = "1.584"
$ws.Range("E50").Value = "  -2.96%  "

# Row 51
$ws.Range("D51").Value = "62.34"
$ws.Range("E51").Value = "  -3.03%  "

# Restore the default (unstyled) look for column D now that the
# values are safely stored as text.
$ws.Range("D2:D51").Style = "Normal"
